# ----------------------------------------------------------------------
# Commit: "Added Invalid code to Login POM Class / Added new sheet to
# excel / Created new InvalidLogin test class / updated IAUtoConst for
# excel sheet"
#
# Data-level changes applied to the workbook:
#   1. ValidLogin sheet: the old single "UserNamePassword" header cell
#      is split into two header cells (UserName / Password) on row 1,
#      keeping the existing admin/manager credential row.
#   2. A new "InvalidLogin" worksheet is added right after ValidLogin,
#      with the same header row and a negative credential pair
#      (Bhanu / Damager).
#   3. Minor cosmetic/theme updates that came along with the resave.
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Theme: swap accent1 <-> accent5 ---
$themeColors = $wb.Theme.ThemeColorScheme
$accent1 = $themeColors.Item(5)
$accent5 = $themeColors.Item(9)
$accent1Rgb = $accent1.RGB
$accent5Rgb = $accent5.RGB
$accent1.RGB = $accent5Rgb
$accent5.RGB = $accent1Rgb

# --- Sheet 1: ValidLogin ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ValidLogin"

$ws1.Range("A1").Value = "UserName"
$ws1.Range("B1").Value = "Password"
$ws1.Range("A2").Value = "admin"
$ws1.Range("B2").Value = "manager"

[void]($ws1.Application.ActiveWindow.Zoom = 235)
$ws1.Range("A1:B2").Select() | Out-Null

# --- Sheet 2: InvalidLogin (new, inserted right after ValidLogin) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "InvalidLogin"

$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "Bhanu"
$ws2.Range("B2").Value = "Damager"

$ws2.Columns.Item(1).AutoFit()
$ws2.Columns.Item(2).AutoFit()

[void]($ws2.Application.ActiveWindow.Zoom = 250)
$ws2.Range("B3").Select() | Out-Null

$ws2.Activate()
